# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2062.25
$ws.Range("J17").Value = 2062.7144
$ws.Range("L17").Value = 6188.1432
$ws.Range("N17").Value = -6524.1432

$ws.Range("H28").Value = 1268.75
$ws.Range("I28").Value = 466.5625
$ws.Range("K28").Value = 466.5625
$ws.Range("M28").Value = 18.4375

$ws.Range("H62").Value = 2635
$ws.Range("J62").Value = 2897
$ws.Range("L62").Value = 2897
$ws.Range("N62").Value = -4145

$ws.Range("H65").Value = 2635
$ws.Range("J65").Value = 2897
$ws.Range("L65").Value = 14485
$ws.Range("N65").Value = -20725

$ws.Range("H96").Value = 414.3
$ws.Range("I96").Value = 368
$ws.Range("J96").Value = 599.5
$ws.Range("K96").Value = 1104
$ws.Range("L96").Value = 1798.5
$ws.Range("M96").Value = 269
$ws.Range("N96").Value = -4544.5

$ws.Range("H98").Value = 1648.6666
$ws.Range("I98").Value = 1639
$ws.Range("K98").Value = 1639
$ws.Range("M98").Value = -141

$ws.Range("H107").Value = 899
$ws.Range("I107").Value = 899
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 899
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1021
$ws.Range("N107").Value = ""

$ws.Range("H113").Value = 3658.5715
$ws.Range("J113").Value = 3880
$ws.Range("L113").Value = 3880
$ws.Range("N113").Value = -10388

$ws.Range("H122").Value = 1648.6666
$ws.Range("I122").Value = 1639
$ws.Range("K122").Value = 4917
$ws.Range("M122").Value = -2467

$ws.Range("H125").Value = 979.1429000000001
$ws.Range("J125").Value = 979.1429000000001
$ws.Range("L125").Value = 8812.286100000001
$ws.Range("N125").Value = -13732.2861

$ws.Range("H129").Value = 2203.0645
$ws.Range("I129").Value = 768.4375
$ws.Range("J129").Value = 3733.3333
$ws.Range("K129").Value = 2305.3125
$ws.Range("L129").Value = 11199.9999
$ws.Range("M129").Value = 2694.6875
$ws.Range("N129").Value = -21199.9999

$ws.Range("H137").Value = 5797.0454
$ws.Range("I137").Value = 1356.7059
$ws.Range("J137").Value = 20894.2
$ws.Range("K137").Value = 4070.1177
$ws.Range("L137").Value = 62682.60000000001
$ws.Range("M137").Value = -1520.1177
$ws.Range("N137").Value = -67782.60000000001

$ws.Range("H138").Value = 300851.06
$ws.Range("J138").Value = 441693.25
$ws.Range("L138").Value = 1325079.75
$ws.Range("N138").Value = -1335359.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3802.625
$ws.Range("I61").Value = 2182.7
$ws.Range("K61").Value = 2182.7
$ws.Range("M61").Value = -1970.7

$ws.Range("H74").Value = 119440.64
$ws.Range("I74").Value = 147115.45
$ws.Range("J74").Value = 2591.4443
$ws.Range("K74").Value = 147115.45
$ws.Range("L74").Value = 2591.4443
$ws.Range("M74").Value = -146241.45
$ws.Range("N74").Value = -4339.4443

$ws.Range("H77").Value = 119440.64
$ws.Range("I77").Value = 147115.45
$ws.Range("J77").Value = 2591.4443
$ws.Range("K77").Value = 735577.25
$ws.Range("L77").Value = 12957.2215
$ws.Range("M77").Value = -731209.25
$ws.Range("N77").Value = -21693.2215

$ws.Range("H81").Value = 50001
$ws.Range("J81").Value = 50001
$ws.Range("L81").Value = 50001
$ws.Range("N81").Value = -51997

$ws.Range("H84").Value = 50001
$ws.Range("J84").Value = 50001
$ws.Range("L84").Value = 150003
$ws.Range("N84").Value = -159987

$ws.Range("H136").Value = 3802.625
$ws.Range("I136").Value = 2182.7
$ws.Range("K136").Value = 6548.099999999999
$ws.Range("M136").Value = -3998.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4143.7144
$ws.Range("I134").Value = 4201.2
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 12603.6
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -10068.6
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 10000000
$ws.Range("I23").Value = 10000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 10000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -9999760
$ws.Range("N23").Value = ""

$ws.Range("H27").Value = 10000000
$ws.Range("I27").Value = 10000000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 10000000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -9999808
$ws.Range("N27").Value = ""

$ws.Range("H31").Value = 5241.0435
$ws.Range("I31").Value = 3534.0588
$ws.Range("K31").Value = 3534.0588
$ws.Range("M31").Value = -3239.0588

$ws.Range("H34").Value = 5241.0435
$ws.Range("I34").Value = 3534.0588
$ws.Range("K34").Value = 3534.0588
$ws.Range("M34").Value = -3332.0588

$ws.Range("H58").Value = 2374.923
$ws.Range("I58").Value = 1636.2593
$ws.Range("J58").Value = 4036.9167
$ws.Range("K58").Value = 1636.2593
$ws.Range("L58").Value = 4036.9167
$ws.Range("M58").Value = -1433.2593
$ws.Range("N58").Value = -4442.9167

$ws.Range("H86").Value = 4650065
$ws.Range("J86").Value = 7747440.5
$ws.Range("L86").Value = 7747440.5
$ws.Range("N86").Value = -7749686.5

$ws.Range("H89").Value = 4650065
$ws.Range("J89").Value = 7747440.5
$ws.Range("L89").Value = 38737202.5
$ws.Range("N89").Value = -38748434.5

$ws.Range("I107").Value = 5556273.5
$ws.Range("J107").Value = 495
$ws.Range("K107").Value = 5556273.5
$ws.Range("L107").Value = 495
$ws.Range("M107").Value = -5554353.5
$ws.Range("N107").Value = -4335

$ws.Range("H132").Value = 3792.913
$ws.Range("I132").Value = 3300.2307
$ws.Range("K132").Value = 9900.6921
$ws.Range("M132").Value = -7370.6921

$ws.Range("H136").Value = 2374.923
$ws.Range("I136").Value = 1636.2593
$ws.Range("J136").Value = 4036.9167
$ws.Range("K136").Value = 4908.7779
$ws.Range("L136").Value = 12110.7501
$ws.Range("M136").Value = -2358.7779
$ws.Range("N136").Value = -17210.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 514.75
$ws.Range("I107").Value = 399.5
$ws.Range("K107").Value = 1198.5
$ws.Range("M107").Value = 721.5

$ws.Range("H123").Value = 2035.4839
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 2605
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 7815
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -12715

$ws.Range("H132").Value = 2334.652
$ws.Range("I132").Value = 1433.8334
$ws.Range("J132").Value = 3317.3635
$ws.Range("K132").Value = 12904.5006
$ws.Range("L132").Value = 29856.2715
$ws.Range("M132").Value = -10374.5006
$ws.Range("N132").Value = -34916.2715

$ws.Range("H141").Value = 9168.714
$ws.Range("I141").Value = 9168.714
$ws.Range("K141").Value = 27506.142
$ws.Range("M141").Value = -22326.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5303.1665
$ws.Range("I113").Value = 5427.5
$ws.Range("J113").Value = 4868
$ws.Range("K113").Value = 5427.5
$ws.Range("L113").Value = 4868
$ws.Range("M113").Value = -3257.5
$ws.Range("N113").Value = -9208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2602.1667
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 3503.25
$ws.Range("K31").Value = 800
$ws.Range("L31").Value = 3503.25
$ws.Range("M31").Value = -552
$ws.Range("N31").Value = -3999.25

$ws.Range("H100").Value = 4075.3157
$ws.Range("I100").Value = 4341.923
$ws.Range("J100").Value = 3497.6667
$ws.Range("K100").Value = 4341.923
$ws.Range("L100").Value = 3497.6667
$ws.Range("M100").Value = -3800.923
$ws.Range("N100").Value = -4579.6667

$ws.Range("H132").Value = 5002.8276
$ws.Range("I132").Value = 3945.7896
$ws.Range("K132").Value = 11837.3688
$ws.Range("M132").Value = -9307.3688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""
